# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Both sheets contain identical data, and the same set of rows changed in each.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 143
    3  = 1669
    4  = 662
    7  = 11752
    11 = 389
    12 = 1103
    13 = 832
    14 = 13421
    15 = 13285
    20 = 263
    23 = 150
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
